# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on the
# zh-cn and de-de worksheets to reflect the new report run.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 12:15:45"
$wsZhCn.Range("H2").Value = "2016-03-17 12:16:04"
$wsZhCn.Range("E4").Value = "2016-03-17 12:15:45"
$wsZhCn.Range("H4").Value = "2016-03-17 12:16:04"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 12:15:48"
$wsDeDe.Range("H2").Value = "2016-03-17 12:16:10"
$wsDeDe.Range("E4").Value = "2016-03-17 12:15:48"
$wsDeDe.Range("H4").Value = "2016-03-17 12:16:10"
